# Update scripts wuth new tpm
# Refresh the LR-pair NATMI TPM-derived metrics (ligand/receptor expression,
# specificity, and edge-weight columns G, H-J, M-T) for every row in the
# Vegfa-Flt1 sheet with the recomputed values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 2.628848666666667
$ws.Range("H2").Value = 7.886546000000001
$ws.Range("I2").Value = 0.04622248078033103
$ws.Range("J2").Value = 0.04850184447997802
$ws.Range("M2").Value = 86.066935
$ws.Range("N2").Value = 258.200805
$ws.Range("O2").Value = 0.8916716774694496
$ws.Range("P2").Value = 0.8955157110805073
$ws.Range("Q2").Value = 226.2569473188367
$ws.Range("R2").Value = 2036.31252586953
$ws.Range("S2").Value = 0.04121527697419717
$ws.Range("T2").Value = 0.0434341637482037

$ws.Range("G3").Value = 2.628848666666667
$ws.Range("H3").Value = 7.886546000000001
$ws.Range("I3").Value = 0.04622248078033103
$ws.Range("J3").Value = 0.04850184447997802
$ws.Range("O3").Value = 0.003820894467605101
$ws.Range("P3").Value = 0.003837366502243974
$ws.Range("Q3").Value = 0.9695316562271111
$ws.Range("R3").Value = 8.725784906044002
$ws.Range("S3").Value = 0.00017661122109255
$ws.Range("T3").Value = 0.0001861193533045145

$ws.Range("G4").Value = 2.628848666666667
$ws.Range("H4").Value = 7.886546000000001
$ws.Range("I4").Value = 0.04622248078033103
$ws.Range("J4").Value = 0.04850184447997802
$ws.Range("M4").Value = 7.669867666666666
$ws.Range("N4").Value = 23.009603
$ws.Range("O4").Value = 0.07946145367329926
$ws.Range("P4").Value = 0.07980401529819077
$ws.Range("Q4").Value = 20.16292138902644
$ws.Range("R4").Value = 181.466292501238
$ws.Range("S4").Value = 0.00367290551519124
$ws.Range("T4").Value = 0.003870641938870636

$ws.Range("G5").Value = 2.628848666666667
$ws.Range("H5").Value = 7.886546000000001
$ws.Range("I5").Value = 0.04622248078033103
$ws.Range("J5").Value = 0.04850184447997802
$ws.Range("M5").Value = 1.242987
$ws.Range("N5").Value = 2.485974
$ws.Range("O5").Value = 0.01287760861197995
$ws.Range("P5").Value = 0.008622083011467191
$ws.Range("Q5").Value = 3.267624717634
$ws.Range("R5").Value = 19.605748305804
$ws.Range("S5").Value = 0.0005952350165638686
$ws.Range("T5").Value = 0.0004181869293156423

$ws.Range("G6").Value = 2.628848666666667
$ws.Range("H6").Value = 7.886546000000001
$ws.Range("I6").Value = 0.04622248078033103
$ws.Range("J6").Value = 0.04850184447997802
$ws.Range("M6").Value = 1.174528666666667
$ws.Range("N6").Value = 3.523586
$ws.Range("O6").Value = 0.01216836577766621
$ws.Range("P6").Value = 0.01222082410759068
$ws.Range("Q6").Value = 3.087658119328444
$ws.Range("R6").Value = 27.788923073956
$ws.Range("S6").Value = 0.0005624520532862145
$ws.Range("T6").Value = 0.0005927325102835294

$ws.Range("I7").Value = 0.472133375270229
$ws.Range("J7").Value = 0.4954156322762335
$ws.Range("M7").Value = 86.066935
$ws.Range("N7").Value = 258.200805
$ws.Range("O7").Value = 0.8916716774694496
$ws.Range("P7").Value = 0.8955157110805073
$ws.Range("Q7").Value = 2311.071461604397
$ws.Range("R7").Value = 20799.64315443957
$ws.Range("S7").Value = 0.4209879587165183
$ws.Range("T7").Value = 0.4436524822182504

$ws.Range("I8").Value = 0.472133375270229
$ws.Range("J8").Value = 0.4954156322762335
$ws.Range("O8").Value = 0.003820894467605101
$ws.Range("P8").Value = 0.003837366502243974
$ws.Range("S8").Value = 0.001803971801541741
$ws.Range("T8").Value = 0.001901091351984837

$ws.Range("I9").Value = 0.472133375270229
$ws.Range("J9").Value = 0.4954156322762335
$ws.Range("M9").Value = 7.669867666666666
$ws.Range("N9").Value = 23.009603
$ws.Range("O9").Value = 0.07946145367329926
$ws.Range("P9").Value = 0.07980401529819077
$ws.Range("Q9").Value = 205.9514757754024
$ws.Range("R9").Value = 1853.563281978622
$ws.Range("S9").Value = 0.03751640432665372
$ws.Range("T9").Value = 0.03953615669713539

$ws.Range("I10").Value = 0.472133375270229
$ws.Range("J10").Value = 0.4954156322762335
$ws.Range("M10").Value = 1.242987
$ws.Range("N10").Value = 2.485974
$ws.Range("O10").Value = 0.01287760861197995
$ws.Range("P10").Value = 0.008622083011467191
$ws.Range("Q10").Value = 33.376717584346
$ws.Range("R10").Value = 200.260305506076
$ws.Range("S10").Value = 0.006079948819383063
$ws.Range("T10").Value = 0.00427151470666419

$ws.Range("I11").Value = 0.472133375270229
$ws.Range("J11").Value = 0.4954156322762335
$ws.Range("M11").Value = 1.174528666666667
$ws.Range("N11").Value = 3.523586
$ws.Range("O11").Value = 0.01216836577766621
$ws.Range("P11").Value = 0.01222082410759068
$ws.Range("Q11").Value = 31.53847272904044
$ws.Range("R11").Value = 283.846254561364
$ws.Range("S11").Value = 0.005745091606132295
$ws.Range("T11").Value = 0.006054387302198675

$ws.Range("G12").Value = 11.96574466666667
$ws.Range("H12").Value = 35.897234
$ws.Range("I12").Value = 0.2103911152781009
$ws.Range("J12").Value = 0.2207661073338543
$ws.Range("M12").Value = 86.066935
$ws.Range("N12").Value = 258.200805
$ws.Range("O12").Value = 0.8916716774694496
$ws.Range("P12").Value = 0.8955157110805073
$ws.Range("Q12").Value = 1029.854968452596
$ws.Range("R12").Value = 9268.694716073369
$ws.Range("S12").Value = 0.1875997986846926
$ws.Range("T12").Value = 0.1976995175915521

$ws.Range("G13").Value = 11.96574466666667
$ws.Range("H13").Value = 35.897234
$ws.Range("I13").Value = 0.2103911152781009
$ws.Range("J13").Value = 0.2207661073338543
$ws.Range("O13").Value = 0.003820894467605101
$ws.Range("P13").Value = 0.003837366502243974
$ws.Range("Q13").Value = 4.413022473208444
$ws.Range("R13").Value = 39.717202258876
$ws.Range("S13").Value = 0.0008038822483993628
$ws.Range("T13").Value = 0.0008471604651137302

$ws.Range("G14").Value = 11.96574466666667
$ws.Range("H14").Value = 35.897234
$ws.Range("I14").Value = 0.2103911152781009
$ws.Range("J14").Value = 0.2207661073338543
$ws.Range("M14").Value = 7.669867666666666
$ws.Range("N14").Value = 23.009603
$ws.Range("O14").Value = 0.07946145367329926
$ws.Range("P14").Value = 0.07980401529819077
$ws.Range("Q14").Value = 91.77567812645576
$ws.Range("R14").Value = 825.9811031381018
$ws.Range("S14").Value = 0.01671798385994458
$ws.Range("T14").Value = 0.01761802180699293

$ws.Range("G15").Value = 11.96574466666667
$ws.Range("H15").Value = 35.897234
$ws.Range("I15").Value = 0.2103911152781009
$ws.Range("J15").Value = 0.2207661073338543
$ws.Range("M15").Value = 1.242987
$ws.Range("N15").Value = 2.485974
$ws.Range("O15").Value = 0.01287760861197995
$ws.Range("P15").Value = 0.008622083011467191
$ws.Range("Q15").Value = 14.873265065986
$ws.Range("R15").Value = 89.23959039591598
$ws.Range("S15").Value = 0.002709334437989338
$ws.Range("T15").Value = 0.001903463703550967

$ws.Range("G16").Value = 11.96574466666667
$ws.Range("H16").Value = 35.897234
$ws.Range("I16").Value = 0.2103911152781009
$ws.Range("J16").Value = 0.2207661073338543
$ws.Range("M16").Value = 1.174528666666667
$ws.Range("N16").Value = 3.523586
$ws.Range("O16").Value = 0.01216836577766621
$ws.Range("P16").Value = 0.01222082410759068
$ws.Range("Q16").Value = 14.05411012901378
$ws.Range("R16").Value = 126.486991161124
$ws.Range("S16").Value = 0.00256011604707507
$ws.Range("T16").Value = 0.002697943766644518

$ws.Range("G17").Value = 8.018423
$ws.Range("H17").Value = 16.036846
$ws.Range("I17").Value = 0.1409862072722574
$ws.Range("J17").Value = 0.09862576223372788
$ws.Range("M17").Value = 86.066935
$ws.Range("N17").Value = 258.200805
$ws.Range("O17").Value = 0.8916716774694496
$ws.Range("P17").Value = 0.8955157110805073
$ws.Range("Q17").Value = 690.1210911435051
$ws.Range("R17").Value = 4140.72654686103
$ws.Range("S17").Value = 0.1257134079385093
$ws.Range("T17").Value = 0.08832091959759386

$ws.Range("G18").Value = 8.018423
$ws.Range("H18").Value = 16.036846
$ws.Range("I18").Value = 0.1409862072722574
$ws.Range("J18").Value = 0.09862576223372788
$ws.Range("O18").Value = 0.003820894467605101
$ws.Range("P18").Value = 0.003837366502243974
$ws.Range("Q18").Value = 2.957231821707333
$ws.Range("R18").Value = 17.743390930244
$ws.Range("S18").Value = 0.0005386934193751945
$ws.Range("T18").Value = 0.0003784631962539862

$ws.Range("G19").Value = 8.018423
$ws.Range("H19").Value = 16.036846
$ws.Range("I19").Value = 0.1409862072722574
$ws.Range("J19").Value = 0.09862576223372788
$ws.Range("M19").Value = 7.669867666666666
$ws.Range("N19").Value = 23.009603
$ws.Range("O19").Value = 0.07946145367329926
$ws.Range("P19").Value = 0.07980401529819077
$ws.Range("Q19").Value = 61.50024330535633
$ws.Range("R19").Value = 369.001459832138
$ws.Range("S19").Value = 0.01120296897773865
$ws.Range("T19").Value = 0.007870731838096145

$ws.Range("G20").Value = 8.018423
$ws.Range("H20").Value = 16.036846
$ws.Range("I20").Value = 0.1409862072722574
$ws.Range("J20").Value = 0.09862576223372788
$ws.Range("M20").Value = 1.242987
$ws.Range("N20").Value = 2.485974
$ws.Range("O20").Value = 0.01287760861197995
$ws.Range("P20").Value = 0.008622083011467191
$ws.Range("Q20").Value = 9.966795549500999
$ws.Range("R20").Value = 39.867182198004
$ws.Range("S20").Value = 0.001815565196939613
$ws.Range("T20").Value = 0.0008503595090484276

$ws.Range("G21").Value = 8.018423
$ws.Range("H21").Value = 16.036846
$ws.Range("I21").Value = 0.1409862072722574
$ws.Range("J21").Value = 0.09862576223372788
$ws.Range("M21").Value = 1.174528666666667
$ws.Range("N21").Value = 3.523586
$ws.Range("O21").Value = 0.01216836577766621
$ws.Range("P21").Value = 0.01222082410759068
$ws.Range("Q21").Value = 9.417867674959332
$ws.Range("R21").Value = 56.507206049756
$ws.Range("S21").Value = 0.001715571739694693
$ws.Range("T21").Value = 0.001205288092735448

$ws.Range("G22").Value = 7.408770666666666
$ws.Range("H22").Value = 22.226312
$ws.Range("I22").Value = 0.1302668213990815
$ws.Range("J22").Value = 0.1366906536762062
$ws.Range("M22").Value = 86.066935
$ws.Range("N22").Value = 258.200805
$ws.Range("O22").Value = 0.8916716774694496
$ws.Range("P22").Value = 0.8955157110805073
$ws.Range("Q22").Value = 637.6501833979066
$ws.Range("R22").Value = 5738.85165058116
$ws.Range("S22").Value = 0.1161552351555322
$ws.Range("T22").Value = 0.1224086279249072

$ws.Range("G23").Value = 7.408770666666666
$ws.Range("H23").Value = 22.226312
$ws.Range("I23").Value = 0.1302668213990815
$ws.Range("J23").Value = 0.1366906536762062
$ws.Range("O23").Value = 0.003820894467605101
$ws.Range("P23").Value = 0.003837366502243974
$ws.Range("Q23").Value = 2.732389196129778
$ws.Range("R23").Value = 24.591502765168
$ws.Range("S23").Value = 0.0004977357771962525
$ws.Range("T23").Value = 0.0005245321355869057

$ws.Range("G24").Value = 7.408770666666666
$ws.Range("H24").Value = 22.226312
$ws.Range("I24").Value = 0.1302668213990815
$ws.Range("J24").Value = 0.1366906536762062
$ws.Range("M24").Value = 7.669867666666666
$ws.Range("N24").Value = 23.009603
$ws.Range("O24").Value = 0.07946145367329926
$ws.Range("P24").Value = 0.07980401529819077
$ws.Range("Q24").Value = 56.8242905860151
$ws.Range("R24").Value = 511.418615274136
$ws.Range("S24").Value = 0.01035119099377107
$ws.Range("T24").Value = 0.01090846301709565

$ws.Range("G25").Value = 7.408770666666666
$ws.Range("H25").Value = 22.226312
$ws.Range("I25").Value = 0.1302668213990815
$ws.Range("J25").Value = 0.1366906536762062
$ws.Range("M25").Value = 1.242987
$ws.Range("N25").Value = 2.485974
$ws.Range("O25").Value = 0.01287760861197995
$ws.Range("P25").Value = 0.008622083011467191
$ws.Range("Q25").Value = 9.209005624647999
$ws.Range("R25").Value = 55.25403374788799
$ws.Range("S25").Value = 0.001677525141104066
$ws.Range("T25").Value = 0.001178558162887963

$ws.Range("G26").Value = 7.408770666666666
$ws.Range("H26").Value = 22.226312
$ws.Range("I26").Value = 0.1302668213990815
$ws.Range("J26").Value = 0.1366906536762062
$ws.Range("M26").Value = 1.174528666666667
$ws.Range("N26").Value = 3.523586
$ws.Range("O26").Value = 0.01216836577766621
$ws.Range("P26").Value = 0.01222082410759068
$ws.Range("Q26").Value = 8.701813532759109
$ws.Range("R26").Value = 78.31632179483199
$ws.Range("S26").Value = 0.001585134331477941
$ws.Range("T26").Value = 0.001670472435728509
